# 自动更新Excel文件 - Mon Dec 15 23:27:04 UTC 2025
#
# Daily rollover: for every data row, the "剩余" (E, days remaining) count
# ticks down by one day. When a row has already hit 1 remaining day, it
# means the container was refilled/serviced that day, so instead of going
# to 0 the row resets: 剩余 goes back up to the full 总天 (D) count and the
# 开始时间 (F) start-date is pushed forward by one full cycle (+7 days,
# matching the existing YYYYMMDD-style integer date encoding).
#
# Row 36 is intentionally skipped - its 开始时间 value is corrupted
# ("202510929", not a real YYYYMMDD date), so it is left untouched rather
# than guessing a rollover.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {

    $totalDays = $ws.Cells.Item($row, 4).Value2   # D: 总天
    $remaining = $ws.Cells.Item($row, 5).Value2   # E: 剩余
    $startDate = $ws.Cells.Item($row, 6).Value2   # F: 开始时间 (YYYYMMDD integer)

    if ($null -eq $remaining) { continue }

    # Skip rows whose start-date isn't a sane 8-digit YYYYMMDD value -
    # e.g. row 36's "202510929" - since the +7 day rollover math would be
    # meaningless on a corrupted date.
    $dateText = ""
    if ($null -ne $startDate) {
        $dateText = [string][int64]$startDate
        if ($dateText.Length -ne 8) {
            continue
        }
    }

    if ($remaining -eq 1) {
        # Refill day: reset remaining to the full cycle length and bump
        # the start date forward by 7 days (still as a YYYYMMDD integer).
        $year  = [int]$dateText.Substring(0, 4)
        $month = [int]$dateText.Substring(4, 2)
        $day   = [int]$dateText.Substring(6, 2)
        $asDate = Get-Date -Year $year -Month $month -Day $day
        $newDate = $asDate.AddDays(7)
        $newDateNum = [int]$newDate.ToString("yyyyMMdd")

        $ws.Cells.Item($row, 5).Value2 = $totalDays
        $ws.Cells.Item($row, 6).Value2 = $newDateNum
    }
    else {
        # Normal day: one fewer day remains before the next refill.
        $ws.Cells.Item($row, 5).Value2 = $remaining - 1
    }
}
